# ind_study_region_matrix.xlsx amendment
# "open_space_defs" sheet: insert two new columns capturing point/line
# attributes used to intersect open spaces, and split out a distinct
# "linear_waterway" tag-value class (separated from the general water
# tags because of its high potential for unintentional aggregation on
# proximity across large areas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("open_space_defs")
$ws.Activate()

# --- Insert the two new columns -------------------------------------------------
# New column I: "linear_waterway" (list of waterway tag values treated as
# linear features). This pushes the old I..O columns one to the right.
$ws.Range("I1").EntireColumn.Insert()

# New column K (after the shift above, this sits right after the old I
# column, which is now J): "point_line_fields" (point/line attribute keys
# used to tag open spaces from intersecting points/lines).
$ws.Range("K1").EntireColumn.Insert()

# --- Column widths for the two newly-inserted columns ----------------------------
# (Closest values reachable through this host's ColumnWidth rounding.)
$ws.Columns("I").ColumnWidth = 15
$ws.Columns("K").ColumnWidth = 33

# --- Header row --------------------------------------------------------------
$ws.Range("I1").Value = "linear_waterway"
$ws.Range("K1").Value = "point_line_fields"

# --- linear_waterway values (column I) ----------------------------------------
$linearWaterway = @{
  2  = "river"
  3  = "riverbank"
  4  = "riverbed"
  5  = "strait"
  6  = "waterway"
  7  = "stream"
  8  = "ditch"
  9  = "river"
  10 = "drain"
  11 = "canal"
  12 = "rapids"
  13 = "drystream"
  14 = "brook"
  15 = "derelict_canal"
  16 = "fairway"
}
foreach ($row in $linearWaterway.Keys) {
  $ws.Cells.Item($row, 9).Value = $linearWaterway[$row]
}

# --- point_line_fields values (column K) ---------------------------------------
$pointLineFields = @{
  2 = "amenity"
  3 = "leisure"
  4 = '"natural"'
  5 = "tourism"
  6 = "waterway"
}
foreach ($row in $pointLineFields.Keys) {
  $ws.Cells.Item($row, 11).Value = $pointLineFields[$row]
}

# --- View state: active cell / scroll position ----------------------------------
$ws.Range("J9").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

Write-Output "open_space_defs: inserted linear_waterway + point_line_fields columns"
